# A new price-quote record was inserted into the "Membrillo" price list
# right after the existing row 158 (before row 159), pushing every
# subsequent record down by one row (old row 159 -> new row 160, ...,
# old row 273 -> new row 274). This mirrors the canonical-XML diff, where
# every row from 159..273 ends up holding the values that used to belong
# to the row immediately above it, and a brand new row 274 is created.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one blank row at position 159; Excel automatically shifts rows
# 159:273 down to 160:274 and extends the used range / dimension.
$ws.Rows("159:159").Insert()

# Populate the newly inserted row 159 with the new record's data.
$ws.Range("A159").Value = 10
$ws.Range("B159").Value = "Vega Modelo de Temuco"
$ws.Range("C159").Value = "La Araucanía"
$ws.Range("D159").Value = 45062
$ws.Range("E159").Value = 9
$ws.Range("F159").Value = "Fruta"
$ws.Range("G159").Value = 100104
$ws.Range("H159").Value = "Frutos de pepita"
$ws.Range("I159").Value = 100104003
$ws.Range("J159").Value = "Membrillo"
$ws.Range("K159").Value = "Champion"
$ws.Range("L159").Value = "Primera"
$ws.Range("M159").Value = 100
$ws.Range("N159").Value = 14000
$ws.Range("O159").Value = 14000
$ws.Range("P159").Value = 14000
$ws.Range("Q159").Value = "$/bandeja 18 kilos granel"
$ws.Range("R159").Value = "Región de O'Higgins"
$ws.Range("S159").Value = 778
$ws.Range("T159").Value = 18
